$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11403
$ws.Range("C3").Value = 10555
$ws.Range("C4").Value = 10555
$ws.Range("C5").Value = 10555
$ws.Range("C6").Value = 10550
$ws.Range("C7").Value = 10127
$ws.Range("C8").Value = 10127
$ws.Range("C9").Value = 9747
$ws.Range("C10").Value = 9747
$ws.Range("C11").Value = 9747
$ws.Range("C12").Value = 9747
$ws.Range("C13").Value = 9738
$ws.Range("C14").Value = 9738
$ws.Range("C15").Value = 9738
$ws.Range("C16").Value = 9738
$ws.Range("C17").Value = 9738
$ws.Range("C18").Value = 9738
$ws.Range("C19").Value = 8873
$ws.Range("C20").Value = 8810
$ws.Range("C21").Value = 8810
$ws.Range("C22").Value = 8810
$ws.Range("C23").Value = 8810
$ws.Range("C24").Value = 8810
$ws.Range("C25").Value = 8810
$ws.Range("C26").Value = 8810
$ws.Range("C27").Value = 8810
$ws.Range("C28").Value = 8810
$ws.Range("C29").Value = 8810
$ws.Range("C30").Value = 7736
$ws.Range("C31").Value = 7736
$ws.Range("C32").Value = 7736
$ws.Range("C33").Value = 7736
$ws.Range("C34").Value = 7736
$ws.Range("C35").Value = 7736
$ws.Range("C36").Value = 7736
$ws.Range("C37").Value = 7736
$ws.Range("C38").Value = 7736
$ws.Range("C39").Value = 7736
$ws.Range("C40").Value = 7736
$ws.Range("C41").Value = 7736
$ws.Range("C42").Value = 7736
$ws.Range("C43").Value = 7736
$ws.Range("C44").Value = 7736
$ws.Range("C45").Value = 7736
$ws.Range("C46").Value = 7736
$ws.Range("C47").Value = 7736
$ws.Range("C48").Value = 7736
$ws.Range("C49").Value = 7736
$ws.Range("C50").Value = 7736
$ws.Range("C51").Value = 7736
$ws.Range("C52").Value = 7736
$ws.Range("C53").Value = 7736
$ws.Range("C54").Value = 7736
$ws.Range("C55").Value = 7736
$ws.Range("C56").Value = 7736
$ws.Range("C57").Value = 7736
$ws.Range("C93").Value = 7295
$ws.Range("C94").Value = 7295
$ws.Range("C95").Value = 7295
$ws.Range("C96").Value = 7295
$ws.Range("C97").Value = 7295
$ws.Range("C98").Value = 7295
$ws.Range("C99").Value = 7295
$ws.Range("C100").Value = 7295
$ws.Range("C101").Value = 7295
$ws.Range("C102").Value = 7295
$ws.Range("C103").Value = 7295
$ws.Range("C104").Value = 7295
$ws.Range("C105").Value = 7295
$ws.Range("C106").Value = 7295
$ws.Range("C107").Value = 7295
$ws.Range("C108").Value = 7295
$ws.Range("C109").Value = 7295
$ws.Range("C110").Value = 7295
$ws.Range("C111").Value = 7295
$ws.Range("C112").Value = 7295
$ws.Range("C113").Value = 7295
$ws.Range("C114").Value = 7295
$ws.Range("C115").Value = 7295
$ws.Range("C116").Value = 7295
$ws.Range("C117").Value = 7295
$ws.Range("C118").Value = 7295
$ws.Range("C119").Value = 7295
$ws.Range("C120").Value = 7295
$ws.Range("C121").Value = 7295
$ws.Range("C122").Value = 7295
$ws.Range("C123").Value = 7295
$ws.Range("C124").Value = 7295
$ws.Range("C125").Value = 7295
$ws.Range("C126").Value = 7295
$ws.Range("C127").Value = 7295
$ws.Range("C128").Value = 7295
$ws.Range("C129").Value = 7295
$ws.Range("C130").Value = 7295
$ws.Range("C131").Value = 7295
$ws.Range("C132").Value = 7295
$ws.Range("C133").Value = 7295
$ws.Range("C134").Value = 7295
$ws.Range("C135").Value = 7295
$ws.Range("C136").Value = 7295
